# Zeitplan.xlsx update: continue the project schedule.
#  - The "Testphase" block (F15:G18) that used to carry its own red shades
#    is merged visually into the same fill used by the "Implementieren"
#    block above it (F11:G14): the three in-between rows (F15-F17, which
#    held "Testphase", "Bugfixes und evtl " and "weitere Tests") are
#    cleared out, and the surviving "Testphase" label moves down to F19
#    (replacing "Fertigstellung der App").
#  - The little two-column legend on the right (I1:J19) is re-colored to
#    match that same shade.
#  - The running note in B20 is replaced: "Allgemeiner Zeitpuffer falls
#    verzögerungen oder Probleme auftauchen" -> "Schriftliche Arbeit ".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the still-needed "Testphase" label before the block is cleared.
$testphase = $ws.Range("F15").Value2

# Re-color F15:G18 to match F11:G14 (format only, keep whatever values are
# still there until we clear them below).
$ws.Range("F11:G14").Copy()
$ws.Range("F15:G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the three superseded labels.
$ws.Range("F15").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("F17").ClearContents()

# "Testphase" now lives on row 19 (replacing "Fertigstellung der App").
$ws.Range("F19").Value2 = $testphase

# Re-color the I:J legend column to the same shade as the schedule bars.
$ws.Range("F11").Copy()
$ws.Range("I1:J19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the running note under the schedule.
$ws.Range("B20").Value2 = "Schriftliche Arbeit "
